$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$syncTimestamp = "2025-08-21 21:17:49"

# Rows 7-26: full update including Enviado (Q) timestamp
$rows = @(
    @{ Row = 7;  Result = "Home Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 8;  Result = "Home Win"; Resultado = "Fallo";   Profit = -3.6;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 9;  Result = "Home Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 10; Result = "Away Win"; Resultado = "Fallo";   Profit = -3.9;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 11; Result = "Home Win"; Resultado = "Fallo";   Profit = -4.8;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 12; Result = "Away Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 13; Result = "Home Win"; Resultado = "Fallo";   Profit = -1.9;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 14; Result = "Home Win"; Resultado = "Fallo";   Profit = -2.7;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 15; Result = "Away Win"; Resultado = "Acierto"; Profit = 16.12; ROI = 375;  Enviado = $syncTimestamp },
    @{ Row = 16; Result = "Away Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 17; Result = "Draw";     Resultado = "Fallo";   Profit = -1.8;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 18; Result = "Away Win"; Resultado = "Fallo";   Profit = -4.2;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 19; Result = "Draw";     Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 20; Result = "Away Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 21; Result = "Home Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 22; Result = "Draw";     Resultado = "Fallo";   Profit = -4.5;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 23; Result = "Home Win"; Resultado = "Fallo";   Profit = -1.4;  ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 24; Result = "Home Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 25; Result = "Home Win"; Resultado = "Fallo";   Profit = -5;    ROI = -100; Enviado = $syncTimestamp },
    @{ Row = 26; Result = "Home Win"; Resultado = "Fallo";   Profit = -3.4;  ROI = -100; Enviado = $syncTimestamp }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 12).Value = "Completed"     # L - Status
    $ws.Cells.Item($row, 13).Value = $r.Result        # M - Result
    $ws.Cells.Item($row, 14).Value = $r.Resultado     # N - Resultado_Real
    $ws.Cells.Item($row, 15).Value = $r.Profit        # O - Profit
    $ws.Cells.Item($row, 16).Value = $r.ROI           # P - ROI
    $ws.Cells.Item($row, 17).Value = $r.Enviado       # Q - Enviado
}

# Rows 27-29: update Status/Result/Resultado_Real/Profit/ROI, but leave Enviado (Q) untouched
$rowsNoEnviado = @(
    @{ Row = 27; Result = "Home Win"; Resultado = "Fallo";   Profit = -5;   ROI = -100 },
    @{ Row = 28; Result = "Away Win"; Resultado = "Acierto"; Profit = 37.5; ROI = 750 },
    @{ Row = 29; Result = "Home Win"; Resultado = "Fallo";   Profit = -5;   ROI = -100 }
)

foreach ($r in $rowsNoEnviado) {
    $row = $r.Row
    $ws.Cells.Item($row, 12).Value = "Completed"     # L - Status
    $ws.Cells.Item($row, 13).Value = $r.Result        # M - Result
    $ws.Cells.Item($row, 14).Value = $r.Resultado     # N - Resultado_Real
    $ws.Cells.Item($row, 15).Value = $r.Profit        # O - Profit
    $ws.Cells.Item($row, 16).Value = $r.ROI           # P - ROI
}
